$wb = $excel.ActiveWorkbook

# Target column width as persisted in the workbook XML (17.2159881591797 -> 13.4101845877511).
# The engine stores ColumnWidth as round(ColumnWidth*6)/6 + 5/6, so 12.5 is the input value
# that lands closest to the target stored width of 13.4101845877511.
$newColumnWidth = 12.5

# Sheet "Overview": columns E (zh-cn) and F (de-de) hold the status text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# Sheet "zh-cn": column C (Status) holds the status text
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# Sheet "de-de": column C (Status) holds the status text
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
